$wb = $excel.ActiveWorkbook

# --- Update "Logs" sheet: append new rows 60-63 ---
$ws = $wb.Worksheets.Item("Logs")

$newRows = @(
    @("Klacht over levering", "mailmind.test@zohomail.eu", "Ik ben niet tevreden over mijn bestelling. Ik hoor graag hoe jullie dit oplossen.", "Klacht", "2025-06-17 22:50:20", "Nee"),
    @("Offerte voor zakelijke samenwerking", "mailmind.test@zohomail.eu", "Kunt u mij een offerte sturen voor 100 stuks product X?", "Bestelling", "2025-06-17 22:56:10", "Nee"),
    @("Offerte voor zakelijke samenwerking", "mailmind.test@zohomail.eu", "Kunt u mij een offerte sturen voor 100 stuks product X?", "Bestelling", "2025-06-17 22:57:10", "Nee"),
    @("Afmelding nieuwsbrief", "mailmind.test@zohomail.eu", "Graag afmelden voor de nieuwsbrief. Dank u.", "Afmelding", "2025-06-17 22:58:10", "Nee")
)

$startRow = 60
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $rowData[0]
    $ws.Cells.Item($r, 2).Value = $rowData[1]
    $ws.Cells.Item($r, 3).Value = $rowData[2]
    $ws.Cells.Item($r, 4).Value = $rowData[3]
    # Column E (Antwoord) intentionally left blank, as in the source data
    $ws.Cells.Item($r, 6).Value = $rowData[4]
    $ws.Cells.Item($r, 7).Value = $rowData[5]
}

# --- Extend conditional formatting ranges to cover the new rows ---
$dFc = $ws.Range("D2:D59").FormatConditions.Item(1)
$dFc.ModifyAppliesToRange($ws.Range("D2:D63"))

$gFc = $ws.Range("G2:G59").FormatConditions.Item(1)
$gFc.ModifyAppliesToRange($ws.Range("G2:G63"))

# --- Update "Dashboard" sheet: refresh category counts ---
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Cells.Item(4, 1).Value = "Afmelding"
$dash.Cells.Item(4, 2).Value = 12

$dash.Cells.Item(5, 1).Value = "Bestelling"
$dash.Cells.Item(5, 2).Value = 6

$dash.Cells.Item(6, 1).Value = "Klacht"
$dash.Cells.Item(6, 2).Value = 6
